# teste_pre_dot.xlsx — add the newly-processed CNPJ "21.578.639/0001-29"
# as a new row at the bottom of the single-column list on Planilha1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new CNPJ value below the existing list (A1:A11 -> A1:A12).
$ws.Range("A12").Value = "21.578.639/0001-29"

# Row 11 (previously the last data row, 18.75pt custom height) goes back
# to the sheet's default row height now that it is no longer the last row.
$ws.Rows(11).AutoFit()

# The new last row (12) picks up the 18.75pt custom height used by the
# other data rows.
$ws.Rows(12).RowHeight = 18.75

# Move the sheet's stored selection from D8 down to the new bottom of the
# list (E10:E11).
$ws.Range("E10:E11").Select()
